$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 170-171; existing rows 170-198 shift down to 172-200.
$ws.Rows("170:171").Insert()

# --- Row 170: new weekly record (Provincia de Cautín, $4000) ---
$ws.Range("A170").Value = 10
$ws.Range("B170").Value = "Vega Modelo de Temuco"
$ws.Range("C170").Value = "La Araucanía"
$ws.Range("D170").Value = 44474
$ws.Range("E170").Value = 9
$ws.Range("F170").Value = 100112044
$ws.Range("G170").Value = "Perejil"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 30
$ws.Range("K170").Value = 4000
$ws.Range("L170").Value = 4000
$ws.Range("M170").Value = 4000
$ws.Range("N170").Value = "$/docena de atados (3 kilos)"
$ws.Range("O170").Value = "Provincia de Cautín"
$ws.Range("P170").Value = 1333
$ws.Range("Q170").Value = 3
$ws.Range("R170").Value = "Hortaliza"

# --- Row 171: new weekly record (Región Metropolitana, $3300) ---
$ws.Range("A171").Value = 10
$ws.Range("B171").Value = "Vega Modelo de Temuco"
$ws.Range("C171").Value = "La Araucanía"
$ws.Range("D171").Value = 44474
$ws.Range("E171").Value = 9
$ws.Range("F171").Value = 100112044
$ws.Range("G171").Value = "Perejil"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 20
$ws.Range("K171").Value = 3300
$ws.Range("L171").Value = 3300
$ws.Range("M171").Value = 3300
$ws.Range("N171").Value = "$/docena de atados (3 kilos)"
$ws.Range("O171").Value = "Región Metropolitana"
$ws.Range("P171").Value = 1100
$ws.Range("Q171").Value = 3

$ws.Range("R171").Value = "Hortaliza"
